# Generate Report for Handoff
# Refresh the handoff report: new report GUID, new content hash for the
# generated .xlf targets, and updated handoff timestamps.

$oldGuid = "b418845a-555e-467d-98cc-21ec515b0289"
$newGuid = "bf54c79d-2bc2-4638-a62f-f73ada39e8cc"

$oldHash = "e1a99f02e3da55052c033171e1bb28fc6d8e6c04"
$newHash = "d45d9a5665ab1adfb22c374dc2c0fe4ed495731a"

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = ($ws.Range("A2").Value() -replace $oldGuid, $newGuid)
$ws.Range("D2").Value = "2016-45-17 12:45:33"

foreach ($h in $ws.Hyperlinks) {
    $h.TextToDisplay = ($h.TextToDisplay -replace $oldGuid, $newGuid)
}

# --- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = ($ws.Range("A2").Value() -replace $oldGuid, $newGuid)
$ws.Range("D2").Value = (($ws.Range("D2").Value() -replace $oldGuid, $newGuid) -replace $oldHash, $newHash)
$ws.Range("E2").Value = "2016-03-17 12:45:30"

foreach ($h in $ws.Hyperlinks) {
    $h.TextToDisplay = (($h.TextToDisplay -replace $oldGuid, $newGuid) -replace $oldHash, $newHash)
}

# --- de-de sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = ($ws.Range("A2").Value() -replace $oldGuid, $newGuid)
$ws.Range("D2").Value = (($ws.Range("D2").Value() -replace $oldGuid, $newGuid) -replace $oldHash, $newHash)
$ws.Range("E2").Value = "2016-03-17 12:45:33"

foreach ($h in $ws.Hyperlinks) {
    $h.TextToDisplay = (($h.TextToDisplay -replace $oldGuid, $newGuid) -replace $oldHash, $newHash)
}
